# "Water rendering attempt 2"
#
# Summary of the change:
#  1. Heading "Implementation & Results" gains a new " & Analysis" suffix,
#     becoming "Implementation & Results & Analysis".
#  2. The "Conclusions and Analysis" heading loses its "Analysis" half,
#     becoming "Conclusions and ".
#  3. The two paragraphs that used to sit directly under that heading
#     ("What have I accomplished? ..." and "Did I hit my targets? Can I
#     simulate an ox-bow lake?") move up into the Results section, right
#     after "Lots of diagrams, examples, code snippets, etc." (which also
#     grows a trailing space run).

$d = $word.ActiveDocument
$CR = [char]13

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -eq $text) {
            return $doc.Paragraphs.Item($i)
        }
    }
    return $null
}

# --- Step 1: delete the two paragraphs below "Conclusions and Analysis"
# (they get re-created further up, right after "Lots of diagrams...").
$conclusionsTarget = "Conclusions and Analysis" + $CR
$conclusionsHeading = Find-ParagraphByText $d $conclusionsTarget
$headingIndex = $conclusionsHeading.Index

$firstOld = $d.Paragraphs.Item($headingIndex + 1)
$secondOld = $d.Paragraphs.Item($headingIndex + 2)
$deleteRange = $d.Range($firstOld.Range.Start, $secondOld.Range.End)
$deleteRange.Delete()

# --- Step 2: truncate "Conclusions and Analysis" -> "Conclusions and ".
$conclusionsHeading.Range.Text = "Conclusions and "

# --- Step 3: "Implementation & Results" -> "Implementation & Results & Analysis"
$implTarget = "Implementation & Results" + $CR
$implHeading = Find-ParagraphByText $d $implTarget
$r = $implHeading.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(" & Analysis")

# --- Step 4: "Lots of diagrams, examples, code snippets, etc." gains a
# trailing space run, then two new paragraphs are inserted right after it.
$diagramsTarget = "Lots of diagrams, examples, code snippets, etc." + $CR
$diagramsPara = Find-ParagraphByText $d $diagramsTarget
$r = $diagramsPara.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(" ")

$diagramsPara.Range.InsertParagraphAfter()
$firstNewIndex = $diagramsPara.Index + 1
$firstNew = $d.Paragraphs.Item($firstNewIndex)
$firstNew.Range.InsertBefore("What have I accomplished? Compare with existing models & simulations in terms of realism & representation. Looking back, would I use a node-based or particle-based simulation?")

$firstNew.Range.InsertParagraphAfter()
$secondNewIndex = $firstNew.Index + 1
$secondNew = $d.Paragraphs.Item($secondNewIndex)
$secondNew.Range.InsertBefore("Did I hit my targets? Can I simulate an ox-bow lake?")

Write-Output "done"
